$wb = $excel.ActiveWorkbook

# --- RenaultKiger sheet: append two new rows ---
$wsKiger = $wb.Worksheets.Item("RenaultKiger")
$wsKiger.Range("A5").Value = "Renault Kiger"
$wsKiger.Range("B5").Value = "₹ 5.45 Lakh"
$wsKiger.Range("A6").Value = "Kiger RXE MT"
$wsKiger.Range("B6").Value = "₹ 5.45 Lakh"

# --- HyundaiVenue sheet: append two new rows ---
$wsVenue = $wb.Worksheets.Item("HyundaiVenue")
$wsVenue.Range("A5").Value = "Hyundai Venue"
$wsVenue.Range("B5").Value = "₹ 6.87 Lakh"
$wsVenue.Range("A6").Value = "Hyundai Venue E 1.2 Petrol"
$wsVenue.Range("B6").Value = "₹ 6.87 Lakh"
